# MV_TestData.xlsx edit:
#  - Insert a new 4th table column "Schadstoffklasse_Title_Soll" (Table1: A1:C4 -> A1:D4)
#  - Swap the existing "Kennzeichen" / "Zulassungsland" columns (B <-> C)
#  - Fill the new column with "Schadstoffklasse" / "Klasa emisji spalin" values
#  - Move the active selection to A4
#  - Re-fit the affected column widths
#  - Set page setup (paper size / orientation) on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add a 4th column to the Table (Table1 grows from A1:C4 to A1:D4) ---
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add()

# --- 2. Swap columns B (Kennzeichen) and C (Zulassungsland) for every data row ---
for ($r = 1; $r -le 4; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    $cVal = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}

# --- 3. Populate the new 4th column (header + data) ---
$ws.Range("D1").Value = "Schadstoffklasse_Title_Soll"
$ws.Range("D2").Value = "Schadstoffklasse"
$ws.Range("D3").Value = "Klasa emisji spalin"
$ws.Range("D4").Value = "Schadstoffklasse"

# --- 4. Move the selection to A4 ---
$ws.Range("A4").Select()

# --- 5. Resize the touched columns to their (re-)fitted widths ---
$ws.Columns.Item(2).ColumnWidth = 15.95
$ws.Columns.Item(3).ColumnWidth = 15.95
$ws.Columns.Item(4).ColumnWidth = 26.95
$ws.Columns.Item(8).ColumnWidth = 13.75

# --- 6. Page setup (paper size / orientation) ---
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

Write-Output "edit applied"
